# Correct a data-entry error in "financial-health-of-residents-data.xlsx":
# the "White areas" (column C) and "Nonwhite areas" (column D) median
# credit-score figures were swapped for every city row on both the
# "City data" and "Peer group data" sheets. Swap the two columns back so
# the values line up with their header labels.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "City data"
$ws2 = $wb.Worksheets.Item(2)   # "Peer group data"

# --- City data sheet: every city data row (skips the 9 category-header
# rows, which have no C/D values, and row 75 which is a live AVERAGE()
# formula that will recompute automatically once its inputs are fixed).
$cityRows = @(
    7,8,9,10,
    12,13,14,15,16,17,18,19,20,21,
    23,24,25,26,27,28,29,
    31,32,33,34,35,36,37,
    39,40,41,42,43,44,
    46,47,48,49,50,51,52,53,
    55,56,57,58,59,60,
    62,63,64,65,66,67,68,69,70,71,
    73,74,
    76
)

foreach ($r in $cityRows) {
    $cVal = $ws1.Cells.Item($r, 3).Value()
    $dVal = $ws1.Cells.Item($r, 4).Value()
    $ws1.Cells.Item($r, 3).Value = $dVal
    $ws1.Cells.Item($r, 4).Value = $cVal
}

# --- Peer group data sheet: every summary row (rows 6-16; all are plain
# cached values, no formulas).
for ($r = 6; $r -le 16; $r++) {
    $cVal = $ws2.Cells.Item($r, 3).Value()
    $dVal = $ws2.Cells.Item($r, 4).Value()
    $ws2.Cells.Item($r, 3).Value = $dVal
    $ws2.Cells.Item($r, 4).Value = $cVal
}

# The group-divider rows on "City data" no longer need to stay merged
# across A:R now that the underlying numbers have been corrected.
$dividerRanges = @(
    "A6:R6", "A11:R11", "A22:R22", "A30:R30", "A38:R38",
    "A45:R45", "A54:R54", "A61:R61", "A72:R72"
)
foreach ($rng in $dividerRanges) {
    $ws1.Range($rng).UnMerge()
}

# Leave the workbook focused on "City data", where the fix was made.
$ws2.Range("D6").Select()
$ws1.Activate()
$ws1.Range("D8").Select()
